$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.41"
$ws.Range("E2").Value = "'0.05%"
$ws.Range("D3").Value = "'31.21"
$ws.Range("E3").Value = "'1.19%"
$ws.Range("E4").Value = "'1.17%"
$ws.Range("D5").Value = "'0.07493"
$ws.Range("E5").Value = "'2.61%"
$ws.Range("E6").Value = "'-1.14%"
$ws.Range("D7").Value = "'7.789"
$ws.Range("E7").Value = "'1.50%"
$ws.Range("D8").Value = "'0.9184"
$ws.Range("E8").Value = "'2.04%"
$ws.Range("D9").Value = "'0.09309"
$ws.Range("E9").Value = "'17.52%"
$ws.Range("D10").Value = "'0.1730"
$ws.Range("E10").Value = "'3.03%"
$ws.Range("D11").Value = "'0.08307"
$ws.Range("E11").Value = "'1.70%"
$ws.Range("D12").Value = "'0.03281"
$ws.Range("E12").Value = "'5.92%"
$ws.Range("D13").Value = "'0.09939"
$ws.Range("E13").Value = "'-1.09%"
$ws.Range("D14").Value = "'0.001498"
$ws.Range("E14").Value = "'0.08%"
$ws.Range("D15").Value = "'0.005713"
$ws.Range("E15").Value = "'-1.61%"
$ws.Range("E16").Value = "'-0.39%"
$ws.Range("D17").Value = "'3.776"
$ws.Range("E17").Value = "'1.58%"
$ws.Range("E18").Value = "'3.76%"
$ws.Range("D19").Value = "'0.3347"
$ws.Range("E19").Value = "'0.71%"
$ws.Range("D21").Value = "'4.095"
$ws.Range("E21").Value = "'1.90%"
$ws.Range("D22").Value = "'0.2099"
$ws.Range("E22").Value = "'0.14%"
$ws.Range("D23").Value = "'0.04536"
$ws.Range("E23").Value = "'0.25%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'0.97%"
$ws.Range("D25").Value = "'0.004305"
$ws.Range("E25").Value = "'-6.94%"
$ws.Range("D26").Value = "'0.0001300"
$ws.Range("E26").Value = "'0.18%"
$ws.Range("D27").Value = "'0.0003395"
$ws.Range("E27").Value = "'0.19%"
$ws.Range("D39").Value = "'0.01622"
$ws.Range("E39").Value = "'2.17%"
$ws.Range("D40").Value = "'0.04578"
$ws.Range("E40").Value = "'3.81%"
$ws.Range("D41").Value = "'0.007494"
$ws.Range("E41").Value = "'2.74%"
$ws.Range("D42").Value = "'0.009843"
$ws.Range("E42").Value = "'14.30%"
$ws.Range("D43").Value = "'0.1360"
$ws.Range("E43").Value = "'3.20%"
$ws.Range("D44").Value = "'0.002220"
$ws.Range("E44").Value = "'11.11%"
$ws.Range("D45").Value = "'0.009784"
$ws.Range("E45").Value = "'3.98%"
$ws.Range("D46").Value = "'0.00006099"
$ws.Range("E46").Value = "'2.96%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.19%"
$ws.Range("D48").Value = "'2.654"
$ws.Range("E48").Value = "'18.45%"
$ws.Range("D49").Value = "'0.002000"
$ws.Range("E49").Value = "'-30.86%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'0.19%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'0.19%"
